{"js": "// CIV-11205 Unit tests covered for location update\n//\n// The order document's recital sentence referenced the court's name a\n// second time (\"... at <<courtName>>.\") where it should instead reference\n// the hearing location in more detail: site name, address and postcode.\n// This updates that sentence, and also removes a stray extra blank\n// paragraph that was left between the recital sentence and the\n// \"THE COURT RECORDS THAT\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that contains the recital sentence about who made\n// the order, when, and at which court.\nconst needle = \"This order is made by\";\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // Replace the trailing \"<<courtName>>\" merge field in that sentence with\n  // the three new merge fields for the hearing location.\n  const found = targetParagraph.search(\"<<courtName>>\", { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    found.items[0].insertText(\n      \"<<siteName>> - <<address>> - <<postcode>>\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// Remove the extra blank paragraph directly following the recital sentence\n// (there were two empty paragraphs in a row; only one should remain).\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet recitalIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    recitalIndex = i;\n    break;\n  }\n}\n\nif (\n  recitalIndex !== -1 &&\n  recitalIndex + 2 < paragraphs.items.length &&\n  paragraphs.items[recitalIndex + 1].text === \"\" &&\n  paragraphs.items[recitalIndex + 2].text === \"\"\n) {\n  paragraphs.items[recitalIndex + 2].delete();\n  await context.sync();\n}\n", "ps1": "# CIV-11205 Unit tests covered for location update\n#\n# The order document's recital sentence referenced the court's name a\n# second time (\"... at <<courtName>>.\") where it should instead reference\n# the hearing location in more detail: site name, address and postcode.\n# This updates that sentence, and also removes a stray extra blank\n# paragraph that was left between the recital sentence and the\n# \"THE COURT RECORDS THAT\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the recital sentence about who made\n# the order, when, and at which court.\n$needle = \"This order is made by\"\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*$needle*\") {\n    $targetIndex = $i\n    break\n  }\n}\n\nif ($targetIndex -ne -1) {\n  $target = $d.Paragraphs.Item($targetIndex)\n\n  # Replace the trailing \"<<courtName>>\" merge field in that sentence with\n  # the three new merge fields for the hearing location.\n  $rng = $target.Range\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = \"<<courtName>>\"\n  $rng.Find.Replacement.Text = \"<<siteName>> - <<address>> - <<postcode>>\"\n  $rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n  # Remove the extra blank paragraph directly following the recital\n  # sentence (there were two empty paragraphs in a row; only one should\n  # remain).\n  $next1 = $d.Paragraphs.Item($targetIndex + 1)\n  $next2 = $d.Paragraphs.Item($targetIndex + 2)\n  if ($next1.Range.Text.Trim() -eq \"\" -and $next2.Range.Text.Trim() -eq \"\") {\n    $next2.Range.Delete()\n  }\n}\n"}
